$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3999.6
$ws.Range("H51").Value = 4153.579
$ws.Range("I51").Value = 4551
$ws.Range("J51").Value = 3040.8
$ws.Range("K51").Value = 4551
$ws.Range("L51").Value = 3040.8
$ws.Range("M51").Value = -4067
$ws.Range("N51").Value = -4008.8
$ws.Range("H81").Value = 199000
$ws.Range("J81").Value = 199000
$ws.Range("L81").Value = 199000
$ws.Range("N81").Value = -200996
$ws.Range("H84").Value = 199000
$ws.Range("J84").Value = 199000
$ws.Range("L84").Value = 597000
$ws.Range("N84").Value = -606984
$ws.Range("H92").Value = 743.8889
$ws.Range("I92").Value = 766.8
$ws.Range("J92").Value = 629.3333
$ws.Range("K92").Value = 766.8
$ws.Range("L92").Value = 629.3333
$ws.Range("M92").Value = 481.2
$ws.Range("N92").Value = -3125.3333
$ws.Range("H101").Value = 902.5
$ws.Range("J101").Value = 1450
$ws.Range("L101").Value = 4350
$ws.Range("N101").Value = -7594
$ws.Range("H103").Value = 1331.6666
$ws.Range("J103").Value = 998
$ws.Range("L103").Value = 2994
$ws.Range("N103").Value = -4166
$ws.Range("H125").Value = 4826
$ws.Range("I125").Value = 325
$ws.Range("K125").Value = 2925
$ws.Range("M125").Value = -465
$ws.Range("H129").Value = 1332.6666
$ws.Range("I129").Value = 1332.6666
$ws.Range("K129").Value = 3997.9998
$ws.Range("M129").Value = 1002.0002
$ws.Range("H138").Value = 2798.202
$ws.Range("I138").Value = 2035.1333
$ws.Range("J138").Value = 2934.4644
$ws.Range("K138").Value = 6105.3999
$ws.Range("L138").Value = 8803.393199999999
$ws.Range("M138").Value = -965.3999000000003
$ws.Range("N138").Value = -19083.3932
$ws.Range("H141").Value = 4421.2085
$ws.Range("I141").Value = 4421.2085
$ws.Range("K141").Value = 13263.6255
$ws.Range("M141").Value = -8083.625499999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3803.6216
$ws.Range("I45").Value = 3240.353
$ws.Range("K45").Value = 3240.353
$ws.Range("M45").Value = -2863.353
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = $null
$ws.Range("N55").Value = $null
$ws.Range("H88").Value = 1473.3
$ws.Range("I88").Value = 1236.25
$ws.Range("J88").Value = 1631.3334
$ws.Range("K88").Value = 1236.25
$ws.Range("L88").Value = 1631.3334
$ws.Range("M88").Value = -830.25
$ws.Range("N88").Value = -2443.3334
$ws.Range("H91").Value = 1473.3
$ws.Range("I91").Value = 1236.25
$ws.Range("J91").Value = 1631.3334
$ws.Range("K91").Value = 1236.25
$ws.Range("L91").Value = 1631.3334
$ws.Range("M91").Value = 167.75
$ws.Range("N91").Value = -4439.3334
$ws.Range("H97").Value = 1291.5883
$ws.Range("I97").Value = 1013.61536
$ws.Range("K97").Value = 1013.61536
$ws.Range("M97").Value = -517.61536
$ws.Range("H132").Value = 3019.42
$ws.Range("I132").Value = 2853.775
$ws.Range("J132").Value = 3682
$ws.Range("K132").Value = 8561.325000000001
$ws.Range("L132").Value = 11046
$ws.Range("M132").Value = -6031.325000000001
$ws.Range("N132").Value = -16106
$ws.Range("H139").Value = 94500
$ws.Range("J139").Value = 94500
$ws.Range("L139").Value = 94500
$ws.Range("N139").Value = -104780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 389.16666
$ws.Range("J22").Value = 623.5
$ws.Range("L22").Value = 623.5
$ws.Range("N22").Value = -969.5
$ws.Range("H86").Value = 3392.5881
$ws.Range("I86").Value = 1599.8572
$ws.Range("J86").Value = 4647.5
$ws.Range("K86").Value = 1599.8572
$ws.Range("L86").Value = 4647.5
$ws.Range("M86").Value = -476.8571999999999
$ws.Range("N86").Value = -6893.5
$ws.Range("H89").Value = 3392.5881
$ws.Range("I89").Value = 1599.8572
$ws.Range("J89").Value = 4647.5
$ws.Range("K89").Value = 7999.286
$ws.Range("L89").Value = 23237.5
$ws.Range("M89").Value = -2383.286
$ws.Range("N89").Value = -34469.5
$ws.Range("H94").Value = 957.1429000000001
$ws.Range("I94").Value = 957.1429000000001
$ws.Range("K94").Value = 957.1429000000001
$ws.Range("M94").Value = -506.1429000000001
$ws.Range("H99").Value = 1768.3334
$ws.Range("I99").Value = 1222
$ws.Range("K99").Value = 1222
$ws.Range("M99").Value = 276
$ws.Range("H105").Value = 1959
$ws.Range("I105").Value = 1703.2858
$ws.Range("K105").Value = 1703.2858
$ws.Range("M105").Value = 43.71419999999989
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = $null
$ws.Range("H133").Value = 120001
$ws.Range("J133").Value = 120001
$ws.Range("L133").Value = 120001
$ws.Range("N133").Value = -130121

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = $null
$ws.Range("H94").Value = 1144.9333
$ws.Range("I94").Value = 622
$ws.Range("J94").Value = 1335.091
$ws.Range("K94").Value = 622
$ws.Range("L94").Value = 1335.091
$ws.Range("M94").Value = -171
$ws.Range("N94").Value = -2237.091
$ws.Range("H99").Value = 3468
$ws.Range("J99").Value = 1900
$ws.Range("L99").Value = 1900
$ws.Range("N99").Value = -4896
$ws.Range("H105").Value = 2529
$ws.Range("I105").Value = 2267.5454
$ws.Range("K105").Value = 2267.5454
$ws.Range("M105").Value = -520.5454
$ws.Range("H126").Value = 3468
$ws.Range("J126").Value = 1900
$ws.Range("L126").Value = 5700
$ws.Range("N126").Value = -10640
$ws.Range("H132").Value = 2706.5
$ws.Range("I132").Value = 906
$ws.Range("K132").Value = 2718
$ws.Range("M132").Value = -188
$ws.Range("H134").Value = 1954.2667
$ws.Range("I134").Value = 1951.2858
$ws.Range("K134").Value = 5853.857400000001
$ws.Range("M134").Value = -3318.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 642
$ws.Range("I16").Value = 370
$ws.Range("K16").Value = 1110
$ws.Range("M16").Value = -937
$ws.Range("H23").Value = 55.8
$ws.Range("J23").Value = 64.75
$ws.Range("L23").Value = 194.25
$ws.Range("N23").Value = -664.25
$ws.Range("H33").Value = 396.1579
$ws.Range("I33").Value = 154.41667
$ws.Range("K33").Value = 926.5000200000001
$ws.Range("M33").Value = -643.5000200000001
$ws.Range("H109").Value = 2217.0625
$ws.Range("I109").Value = 988.9091
$ws.Range("J109").Value = 4919
$ws.Range("K109").Value = 2966.7273
$ws.Range("L109").Value = 14757
$ws.Range("M109").Value = -1926.7273
$ws.Range("N109").Value = -16837
$ws.Range("H121").Value = 8603.385
$ws.Range("I121").Value = 426.5
$ws.Range("J121").Value = 15612.143
$ws.Range("K121").Value = 1279.5
$ws.Range("L121").Value = 46836.429
$ws.Range("M121").Value = 30.5
$ws.Range("N121").Value = -49456.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 931.4
$ws.Range("I97").Value = 510.83334
$ws.Range("K97").Value = 510.83334
$ws.Range("M97").Value = -14.83334000000002
$ws.Range("H102").Value = 1384.7667
$ws.Range("J102").Value = 1810.5
$ws.Range("L102").Value = 1810.5
$ws.Range("N102").Value = -5054.5
$ws.Range("H113").Value = 27267.215
$ws.Range("I113").Value = 12664.389
$ws.Range("K113").Value = 12664.389
$ws.Range("M113").Value = -10494.389

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3481.1904
$ws.Range("J7").Value = 3713.8
$ws.Range("L7").Value = 3713.8
$ws.Range("N7").Value = -3937.8
$ws.Range("H16").Value = 1689.091
$ws.Range("I16").Value = 1587.3684
$ws.Range("K16").Value = 1587.3684
$ws.Range("M16").Value = -1417.3684
$ws.Range("H22").Value = 1944
$ws.Range("I22").Value = 1496.5385
$ws.Range("K22").Value = 1496.5385
$ws.Range("M22").Value = -1201.5385
$ws.Range("H27").Value = 1944
$ws.Range("I27").Value = 1496.5385
$ws.Range("K27").Value = 1496.5385
$ws.Range("M27").Value = -1389.5385
$ws.Range("H40").Value = 2194.6667
$ws.Range("I40").Value = 4874.6665
$ws.Range("J40").Value = 1658.6666
$ws.Range("K40").Value = 4874.6665
$ws.Range("L40").Value = 1658.6666
$ws.Range("M40").Value = -4738.6665
$ws.Range("N40").Value = -1930.6666
$ws.Range("H46").Value = 3667.4546
$ws.Range("I46").Value = 1849.5
$ws.Range("J46").Value = 3849.25
$ws.Range("K46").Value = 1849.5
$ws.Range("L46").Value = 3849.25
$ws.Range("M46").Value = -1661.5
$ws.Range("N46").Value = -4225.25
$ws.Range("H55").Value = 2261.7144
$ws.Range("I55").Value = 3295.6
$ws.Range("K55").Value = 3295.6
$ws.Range("M55").Value = -3122.6
$ws.Range("H126").Value = 3481.1904
$ws.Range("J126").Value = 3713.8
$ws.Range("L126").Value = 11141.4
$ws.Range("N126").Value = -16081.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 30026
$ws.Range("I32").Value = 30026
$ws.Range("K32").Value = 30026
$ws.Range("M32").Value = -29709
$ws.Range("H34").Value = 9026
$ws.Range("I34").Value = 9026
$ws.Range("K34").Value = 9026
$ws.Range("M34").Value = -8823

